# Sprint 43 - Day 9 Test Case Summary: fill in the totals that were left
# blank (Total testcase Written / Total Execution / Total Review) and move
# the active selection to the last cell touched (C53).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C51").Value = 7068
$ws.Range("C52").Value = 2510
$ws.Range("C53").Value = 2510

$ws.Range("C53").Select()
